$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values - use leading apostrophe to force text type
# since the original cells are stored as text (inlineStr), not numbers.
$ws.Range("D2").Value = "'250.86"
$ws.Range("D3").Value = "'22.99"
$ws.Range("D5").Value = "'0.05671"
$ws.Range("D6").Value = "'3.411"
$ws.Range("D7").Value = "'6.375"
$ws.Range("D8").Value = "'0.8134"
$ws.Range("D9").Value = "'0.9228"
$ws.Range("D10").Value = "'0.1441"
$ws.Range("D11").Value = "'0.07455"
$ws.Range("D12").Value = "'0.03121"
$ws.Range("D13").Value = "'0.03064"
$ws.Range("D14").Value = "'0.09356"
$ws.Range("D15").Value = "'3.725"
$ws.Range("D16").Value = "'0.001591"
$ws.Range("D17").Value = "'0.04757"
$ws.Range("D18").Value = "'0.0005789"
$ws.Range("D19").Value = "'0.006381"
$ws.Range("D20").Value = "'0.005042"
$ws.Range("D21").Value = "'0.001031"
$ws.Range("D23").Value = "'3.701"
$ws.Range("D24").Value = "'2.181"
$ws.Range("D26").Value = "'0.1389"
$ws.Range("D28").Value = "'0.0002999"
$ws.Range("D41").Value = "'0.006761"
$ws.Range("D44").Value = "'0.008022"
$ws.Range("D45").Value = "'0.00005802"
$ws.Range("D47").Value = "'0.4999"

# Update Volume(1h) text (column E) for rows where the label changed
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
